$wb = $excel.ActiveWorkbook

# --- Update the "Status" text (shared across Overview/zh-cn/de-de via shared string) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("I2").Value = "d4914487-b880-4503-a38d-9885374630c7.md"
$wsZh.Range("J2").Value = "d4914487-b880-4503-a38d-9885374630c7.096ec810ce57c57ac8a000c72e20162b75c6e289.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-11-08 23:32:06"

# --- de-de sheet: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("I2").Value = "d4914487-b880-4503-a38d-9885374630c7.md"
$wsDe.Range("J2").Value = "d4914487-b880-4503-a38d-9885374630c7.096ec810ce57c57ac8a000c72e20162b75c6e289.de-de.xlf"
$wsDe.Range("K2").Value = "2016-11-08 23:32:24"
